# Update the measured "przypisania"/"porownania" counts (columns G/H) for the
# modified Quick Sort variant on rows 2-41 of quick_testyKlas. Column N
# ("SumMod") is a SUM(G:H) formula and recalculates automatically; the chart
# on the sheet reads straight from N2:N41 so it follows along too.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quick_testyKlas")

$data = @(
    @(2, 940625, 990320),
    @(3, 1895511, 2029508),
    @(4, 2828932, 3187054),
    @(5, 3830703, 4549176),
    @(6, 5107226, 6084360),
    @(7, 6238901, 7514827),
    @(8, 7239661, 9145456),
    @(9, 8942810, 11323266),
    @(10, 9431490, 12671432),
    @(11, 11691806, 15038477),
    @(12, 12527370, 16588047),
    @(13, 14479240, 19027607),
    @(14, 15725580, 21749021),
    @(15, 16177867, 23101046),
    @(16, 18265559, 25747753),
    @(17, 19623979, 27889815),
    @(18, 20682123, 30590036),
    @(19, 21362086, 32641962),
    @(20, 22667913, 35741688),
    @(21, 27200803, 39934480),
    @(22, 26206265, 41220365),
    @(23, 29652915, 45258258),
    @(24, 31237419, 48650006),
    @(25, 31854566, 50620629),
    @(26, 34537990, 54732508),
    @(27, 35798192, 58215365),
    @(28, 39382465, 62603484),
    @(29, 38033011, 64120789),
    @(30, 42054907, 68665467),
    @(31, 42797005, 72898932),
    @(32, 44011903, 76003025),
    @(33, 46858856, 80160440),
    @(34, 48937093, 85026313),
    @(35, 50428823, 88389041),
    @(36, 53898727, 93688611),
    @(37, 55139263, 97228870),
    @(38, 55829868, 100748836),
    @(39, 57625595, 105644874),
    @(40, 61549820, 112217880),
    @(41, 64328131, 115697648)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 7).Value = $row[1]   # Column G
    $ws.Cells.Item($r, 8).Value = $row[2]   # Column H
}

# Matches the saved cursor position recorded in the edited workbook.
$ws.Range("S36").Select() | Out-Null
